$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98, shifting existing rows 98-103 down to 99-104.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price entry.
$ws.Cells.Item(98, 1).Value = 2
$ws.Cells.Item(98, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(98, 3).Value = "Coquimbo"
$ws.Cells.Item(98, 4).Value = 45147
$ws.Cells.Item(98, 5).Value = 4
$ws.Cells.Item(98, 6).Value = 100112026
$ws.Cells.Item(98, 7).Value = "Haba"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 2400
$ws.Cells.Item(98, 11).Value = 10000
$ws.Cells.Item(98, 12).Value = 11000
$ws.Cells.Item(98, 13).Value = 10500
$ws.Cells.Item(98, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(98, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(98, 16).Value = 420
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
